$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.381.02'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.863.59'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.017'
$ws.Range('E4').Value = '  +1.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.42'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.015'
$ws.Range('E6').Value = '  +1.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5118'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3912'
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08297'
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.114'
$ws.Range('E10').Value = '  -0.30%  '
$ws.Range('B11').Value = 'Polkadot'
$ws.Range('C11').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.230'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.864.91'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.27'
$ws.Range('E13').Value = '  -2.39%  '
$ws.Range('B14').Value = 'BinanceUSD'
$ws.Range('C14').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.018'
$ws.Range('E14').Value = '  +1.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.184'
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001100'
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.18'
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06710'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.60'
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.016'
$ws.Range('E20').Value = '  +1.34%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.943'
$ws.Range('E21').Value = '  -1.86%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.411.42'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.09'
$ws.Range('E23').Value = '  -0.93%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.253'
$ws.Range('E24').Value = '  -0.36%  '
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.066.66'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.75'
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.65'
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.407'
$ws.Range('E28').Value = '  -3.94%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '126.61'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.1050'
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.032'
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.807'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.641'
$ws.Range('E33').Value = '  +1.10%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.02434'
$ws.Range('E34').Value = '  -0.89%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.167'
$ws.Range('E35').Value = '  -5.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06470'
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2169'
$ws.Range('E37').Value = '  -1.36%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.252'
$ws.Range('E38').Value = '  +1.73%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.182'
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6411'
$ws.Range('E40').Value = '  -1.78%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.957'
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.10'
$ws.Range('E42').Value = '  -1.96%  '
$ws.Range('B43').Value = 'Decentraland'
$ws.Range('C43').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5992'
$ws.Range('E43').Value = '  -2.16%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.90'
$ws.Range('E44').Value = '  -1.77%  '
$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.704'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.286'
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.986'
$ws.Range('E47').Value = '  -1.58%  '
$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.201'
$ws.Range('E48').Value = '  -2.44%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '120.42'
$ws.Range('E49').Value = '  -1.01%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06864'
$ws.Range('E50').Value = '  -0.59%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '76.25'
$ws.Range('E51').Value = '  -2.53%  '
